$data = @(
    ,("aryl_halide_200epochs", "ranking_test1", 5, 0.001, "mean", "sum", 251.051, 0.698, 15.343, 0.574, 15.087)
    ,("aryl_halide_200epochs", "ranking_test1", 4, 0.01, "max", "sum", 268.523, 0.681, 15.671, 0.535, 15.663)
    ,("aryl_halide_200epochs", "ranking_test1", 4, 0.001, "max", "sum", 262.868, 0.6840000000000001, 15.706, 0.525, 15.846)
    ,("aryl_halide_200epochs", "ranking_test1", 5, 0.01, "max", "sum", 257.144, 0.695, 15.468, 0.591, 14.761)
    ,("aryl_halide_200epochs", "ranking_test1", 5, 0.001, "max", "sum", 252.722, 0.698, 15.341, 0.553, 15.436)
    ,("aryl_halide_200epochs", "ranking_test1", 4, 0.001, "mean", "sum", 268.598, 0.6840000000000001, 15.685, 0.524, 15.984)
    ,("aryl_halide_200epochs", "ranking_test1", 3, 0.01, "max", "sum", 251.058, 0.6899999999999999, 15.474, 0.584, 14.898)
    ,("aryl_halide_200epochs", "ranking_test1", 5, 0.01, "max", "mean", 273.913, 0.67, 15.936, 0.525, 15.817)
    ,("aryl_halide_200epochs", "ranking_test1", 4, 0.01, "max", "product", 246.164, 0.716, 14.78, 0.406, 17.4)
    ,("aryl_halide_200epochs", "ranking_test1", 5, 0.001, "mean", "product", 185.978, 0.779, 13.077, 0.639, 13.818)
    ,("aryl_halide_200epochs", "ranking_test1", 4, 0.01, "mean", "mean", 257.569, 0.6919999999999999, 15.465, 0.521, 16.086)
    ,("aryl_halide_200epochs", "ranking_test1", 5, 0.0001, "mean", "sum", 304.533, 0.634, 16.967, 0.338, 18.869)
    ,("aryl_halide_200epochs", "ranking_test1", 4, 0.0001, "max", "sum", 303.137, 0.624, 17.04, 0.417, 17.618)
    ,("aryl_halide_200epochs", "ranking_test1", 5, 0.0001, "max", "sum", 299.845, 0.637, 16.857, 0.428, 17.548)
    ,("aryl_halide_200epochs", "ranking_test1", 4, 0.0001, "mean", "sum", 327.618, 0.602, 17.659, 0.319, 19.115)
    ,("aryl_halide_200epochs", "ranking_test1", 3, 0.001, "max", "sum", 262.354, 0.674, 15.886, 0.518, 15.947)
    ,("aryl_halide_200epochs", "ranking_test1", 3, 0.001, "max", "sum", 270.111, 0.669, 15.986, 0.496, 16.404)
    ,("aryl_halide_200epochs", "ranking_test1", 5, 0.001, "max", "mean", 254.068, 0.694, 15.429, 0.497, 16.465)
    ,("aryl_halide_200epochs", "ranking_test1", 5, 0.0001, "max", "mean", 307.384, 0.617, 17.267, 0.378, 18.218)
    ,("aryl_halide_200epochs", "ranking_test1", 4, 0.001, "max", "product", 189.158, 0.777, 13.117, 0.644, 13.637)
    ,("aryl_halide_200epochs", "ranking_test1", 4, 0.0001, "max", "product", 241.99, 0.731, 14.499, 0.458, 17.051)
    ,("aryl_halide_200epochs", "ranking_test1", 5, 0.0001, "mean", "product", 228.132, 0.723, 14.632, 0.554, 15.314)
    ,("aryl_halide_200epochs", "ranking_test1", 4, 0.001, "mean", "mean", 281.859, 0.67, 16.002, 0.501, 16.419)
    ,("aryl_halide_200epochs", "ranking_test1", 4, 0.0001, "mean", "mean", 352.875, 0.5590000000000001, 18.54, 0.327, 18.76)
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 74
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    for ($c = 0; $c -lt $rec.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $rec[$c]
    }
}
